$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new rows 19-29 (summary / annotation rows) first, so new
#     shared-string entries are created in the same order the author typed
#     them (new block before going back to tweak the earlier rows). ---
$newRows = @(
    @(1, "Ayoub ", 2),
    @(2, "Yassine", 1),
    @(3, "Yassine", 2),
    @(4, "Amine", 1),
    @(5, "Yanis", 1),
    @(6, "Carlos", 2),
    @(7, "Salim", 2),
    @(8, "Rayane", 2),
    @(9, "Sebastien", 1),
    @(12, "Quentin", 1),
    @(13, "Julien", 1)
)

$r = 19
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Update existing rows with corrected player / goal values ---
$ws.Range("B3").Value = "Salim"
$ws.Range("C3").Value = 1

$ws.Range("B5").Value = "Rayane"

$ws.Range("C7").Value = 2

$ws.Range("B9").Value = "Julien"
$ws.Range("C9").Value = 2

$ws.Range("C10").Value = 2

$ws.Range("C12").Value = 1

$ws.Range("C17").Value = 3

# --- Update the active selection to mirror the authored state ---
$ws.Range("F22").Select()
